$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet had a stray header-only row 6 ("grandes regiões e unidades da
# federação") with no data, which pushed every region's numeric data one
# row below its label (e.g. "norte" in row 6 had no data while "rondônia"
# in row 7 showed norte's numbers, and so on down to the last region,
# whose data lived in an extra trailing row 38).
#
# Deleting that empty header row shifts all the data cells up by one row
# so each region label lines up with its own data; the last (now unused)
# row disappears and the orphaned shared string is dropped automatically.
$ws.Rows("6").Delete()
